$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3:H18").Value = 1
